# Added potions Inventory icons
#
# The asset list line for "UI_INV_Potion_001" needs to be highlighted in
# green (RGB 00B050) and the following "UI_INV_Quest_001" line needs to
# become its own directly-formatted run (explicit "automatic" color) so
# that it is no longer merged with the preceding green run. This mirrors
# splitting a single <w:r> that spans the Scythe/Potion/Quest lines into
# three separate runs.

$d = $word.ActiveDocument

# Locate the start of the "UI_INV_Potion_001" line.
$potionFind = $d.Content
$potionFind.Find.Execute("UI_INV_Potion_001", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$potionStart = $potionFind.Start

# Locate the start of the "UI_INV_Quest_001" line (== end of the Potion line).
$questFind = $d.Content
$questFind.Find.Execute("UI_INV_Quest_001", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$questStart = $questFind.Start

# Locate the start of the "UI_INV_Empty_001" line (== end of the Quest line).
$emptyFind = $d.Content
$emptyFind.Find.Execute("UI_INV_Empty_001", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$emptyStart = $emptyFind.Start

# Whole "UI_INV_Potion_001 ... 1hr" line (incl. trailing line break) turns green.
$potionRange = $d.Range($potionStart, $questStart)
$potionRange.Font.Color = 5287936

# Whole "UI_INV_Quest_001 ... 1hr" line (incl. trailing line break) becomes its
# own explicitly-formatted run again (automatic color), splitting it off from
# the newly green-colored Potion run.
$questRange = $d.Range($questStart, $emptyStart)
$questRange.Font.Color = 255
$questRange.Font.Color = -16777216
